$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 74.5
$ws.Range("C2").Value = -1624.73
$ws.Range("D2").Value = 20520.68
$ws.Range("E2").Value = $null
$ws.Range("H2").Value = 8.44
$ws.Range("I2").Value = 12.13
$ws.Range("J2").Value = 16.77
$ws.Range("K2").Value = "12/25-01/25"

# Row 3
$ws.Range("B3").Value = 78.9
$ws.Range("C3").Value = -54.73
$ws.Range("D3").Value = 12417.15
$ws.Range("E3").Value = $null
$ws.Range("G3").Value = 1.24
$ws.Range("H3").Value = 8.17
$ws.Range("I3").Value = 12.08
$ws.Range("J3").Value = 16.31
$ws.Range("K3").Value = "12/25-01/25"

# Row 4
$ws.Range("A4").Value = 0.12
$ws.Range("B4").Value = 77
$ws.Range("C4").Value = 826.23
$ws.Range("D4").Value = 21994.31
$ws.Range("E4").Value = $null
$ws.Range("G4").Value = 1.33
$ws.Range("H4").Value = 8.59
$ws.Range("I4").Value = 12.74
$ws.Range("J4").Value = 15.75
$ws.Range("K4").Value = "12/25-01/25"

# Row 5
$ws.Range("A5").Value = 1.08
$ws.Range("B5").Value = 80
$ws.Range("C5").Value = 1275.41
$ws.Range("D5").Value = 31207.28
$ws.Range("E5").Value = $null
$ws.Range("G5").Value = 1.33
$ws.Range("H5").Value = 8.9
$ws.Range("I5").Value = 12.97
$ws.Range("J5").Value = 15.19
$ws.Range("K5").Value = "12/25-01/25"
$ws.Range("L5").Value = "מיטב השתלמות כללי"
$ws.Range("M5").Value = 880

# Row 6
$ws.Range("A6").Value = 1.32
$ws.Range("B6").Value = 91
$ws.Range("C6").Value = 6564.77
$ws.Range("D6").Value = 19456.81
$ws.Range("E6").Value = $null
$ws.Range("G6").Value = 1.31
$ws.Range("H6").Value = 9.41
$ws.Range("I6").Value = 13.9
$ws.Range("J6").Value = 15.1
$ws.Range("K6").Value = "12/25-01/25"
$ws.Range("L6").Value = "אנליסט השתלמות כללי "
$ws.Range("M6").Value = 962

# Row 7
$ws.Range("A7").Value = 0.48
$ws.Range("B7").Value = 76.9
$ws.Range("C7").Value = -510.31
$ws.Range("D7").Value = 19494.21
$ws.Range("E7").Value = $null
$ws.Range("G7").Value = 1.27
$ws.Range("H7").Value = 8.34
$ws.Range("I7").Value = 11.9
$ws.Range("J7").Value = 14.82
$ws.Range("K7").Value = "12/25-01/25"
$ws.Range("L7").Value = "מגדל השתלמות כללי"
$ws.Range("M7").Value = 579

# Row 8
$ws.Range("A8").Value = 1.32
$ws.Range("B8").Value = 77.3
$ws.Range("C8").Value = 2499.93
$ws.Range("D8").Value = 25747.29
$ws.Range("E8").Value = $null
$ws.Range("G8").Value = 1.42
$ws.Range("H8").Value = 8.92
$ws.Range("I8").Value = 12.49
$ws.Range("J8").Value = 14.73
$ws.Range("K8").Value = "12/25-01/25"
$ws.Range("L8").Value = "מור השתלמות - כללי "
$ws.Range("M8").Value = 12535

# Row 9
$ws.Range("B9").Value = 76.3
$ws.Range("C9").Value = -995.2
$ws.Range("D9").Value = 14558.74
$ws.Range("E9").Value = $null
$ws.Range("H9").Value = 7.8
$ws.Range("I9").Value = 11.16
$ws.Range("J9").Value = 14.47
$ws.Range("K9").Value = "12/25-01/25"

# Row 10
$ws.Range("A10").Value = 0.84
$ws.Range("B10").Value = 88.8
$ws.Range("C10").Value = -1694.23
$ws.Range("D10").Value = 24801.48
$ws.Range("E10").Value = $null
$ws.Range("G10").Value = 1.16
$ws.Range("H10").Value = 8.46
$ws.Range("I10").Value = 12.79
$ws.Range("J10").Value = 13.38
$ws.Range("K10").Value = "12/25-01/25"

# Row 11
$ws.Range("A11").Value = -2.28
$ws.Range("B11").Value = 81.9
$ws.Range("C11").Value = -7924.24
$ws.Range("D11").Value = 32442.42
$ws.Range("E11").Value = $null
$ws.Range("G11").Value = 0.81
$ws.Range("H11").Value = 6.36
$ws.Range("I11").Value = 11.53
$ws.Range("J11").Value = 12.87
$ws.Range("K11").Value = "12/25-01/25"

# Row 12
$ws.Range("B12").Value = "80.5%"
$ws.Range("C12").Value = "-1,637.1"
$ws.Range("D12").Value = 222640.37
$ws.Range("E12").Value = $null
$ws.Range("H12").Value = "7.77%"
$ws.Range("I12").Value = "12.24%"
$ws.Range("J12").Value = "14.71%"
$ws.Range("K12").Value = "12/25-01/25"

# Row 16
$ws.Range("A16").Value = -1.8
$ws.Range("B16").Value = 91.3
$ws.Range("C16").Value = 667.32
$ws.Range("D16").Value = 1530.23
$ws.Range("E16").Value = $null
$ws.Range("G16").Value = 1.18
$ws.Range("H16").Value = 13.04
$ws.Range("I16").Value = 20.67
$ws.Range("J16").Value = 27.94
$ws.Range("K16").Value = "12/25-01/25"

# Row 17
$ws.Range("A17").Value = -3.6
$ws.Range("B17").Value = 90.5
$ws.Range("C17").Value = 578.15
$ws.Range("D17").Value = 2604.26
$ws.Range("E17").Value = $null
$ws.Range("H17").Value = 11.97
$ws.Range("I17").Value = 19.29
$ws.Range("J17").Value = 27.45
$ws.Range("K17").Value = "12/25-01/25"

# Row 18
$ws.Range("A18").Value = -0.96
$ws.Range("B18").Value = 90
$ws.Range("C18").Value = 1673.84
$ws.Range("D18").Value = 4392.05
$ws.Range("E18").Value = $null
$ws.Range("G18").Value = 1.21
$ws.Range("H18").Value = 13.07
$ws.Range("I18").Value = 21.53
$ws.Range("J18").Value = 26.46
$ws.Range("K18").Value = "12/25-01/25"
$ws.Range("L18").Value = "מיטב השתלמות מניות"
$ws.Range("M18").Value = 883

# Row 19
$ws.Range("A19").Value = -1.44
$ws.Range("B19").Value = 90.5
$ws.Range("C19").Value = 591.12
$ws.Range("D19").Value = 2823.83
$ws.Range("E19").Value = $null
$ws.Range("G19").Value = 1.18
$ws.Range("H19").Value = 12.98
$ws.Range("I19").Value = 19.88
$ws.Range("J19").Value = 26.4
$ws.Range("K19").Value = "12/25-01/25"
$ws.Range("L19").Value = "הראל השתלמות מסלול מניות"
$ws.Range("M19").Value = 763

# Row 20
$ws.Range("A20").Value = -2.16
$ws.Range("B20").Value = 89.7
$ws.Range("C20").Value = 1900.05
$ws.Range("D20").Value = 4537.33
$ws.Range("E20").Value = $null
$ws.Range("G20").Value = 1.19
$ws.Range("H20").Value = 12.97
$ws.Range("I20").Value = 20.43
$ws.Range("J20").Value = 25.89
$ws.Range("K20").Value = "12/25-01/25"
$ws.Range("L20").Value = "הפניקס השתלמות מניות"
$ws.Range("M20").Value = 968

# Row 21
$ws.Range("A21").Value = 0.72
$ws.Range("B21").Value = 85.4
$ws.Range("C21").Value = 3652.58
$ws.Range("D21").Value = 11646.4
$ws.Range("E21").Value = $null
$ws.Range("G21").Value = 1.33
$ws.Range("H21").Value = 13.43
$ws.Range("I21").Value = 20.3
$ws.Range("J21").Value = 25.86
$ws.Range("K21").Value = "12/25-01/25"
$ws.Range("L21").Value = "מור השתלמות - מניות "
$ws.Range("M21").Value = 12536

# Row 22
$ws.Range("A22").Value = -1.2
$ws.Range("B22").Value = 93
$ws.Range("C22").Value = 953.94
$ws.Range("D22").Value = 2637.69
$ws.Range("E22").Value = $null
$ws.Range("G22").Value = 1.22
$ws.Range("H22").Value = 13.13
$ws.Range("I22").Value = 20.79
$ws.Range("J22").Value = 25.48
$ws.Range("K22").Value = "12/25-01/25"
$ws.Range("L22").Value = "מגדל השתלמות מניות"
$ws.Range("M22").Value = 869

# Row 23
$ws.Range("A23").Value = 1.08
$ws.Range("B23").Value = 93.3
$ws.Range("C23").Value = 1094.27
$ws.Range("D23").Value = 12845.38
$ws.Range("E23").Value = $null
$ws.Range("G23").Value = 1.25
$ws.Range("H23").Value = 14.07
$ws.Range("I23").Value = 21.59
$ws.Range("J23").Value = 21.84
$ws.Range("K23").Value = "12/25-01/25"

# Row 24
$ws.Range("A24").Value = -4.92
$ws.Range("B24").Value = 88.8
$ws.Range("C24").Value = -822.42
$ws.Range("D24").Value = 6884.6
$ws.Range("E24").Value = $null
$ws.Range("G24").Value = 0.86
$ws.Range("H24").Value = 9.73
$ws.Range("I24").Value = 19.36
$ws.Range("J24").Value = 21.02
$ws.Range("K24").Value = "12/25-01/25"

# Row 25
$ws.Range("A25").Value = 0.48
$ws.Range("B25").Value = 95.2
$ws.Range("C25").Value = 3416.27
$ws.Range("D25").Value = 15970.81
$ws.Range("E25").Value = $null
$ws.Range("G25").Value = 1.13
$ws.Range("H25").Value = 13.54
$ws.Range("I25").Value = 22.38
$ws.Range("J25").Value = 20.74
$ws.Range("K25").Value = "12/25-01/25"

# Row 26
$ws.Range("B26").Value = "91.1%"
$ws.Range("C26").Value = "13,705.1"
$ws.Range("D26").Value = 65872.58
$ws.Range("E26").Value = $null
$ws.Range("H26").Value = "11.98%"
$ws.Range("I26").Value = "20.59%"
$ws.Range("J26").Value = "23.19%"
$ws.Range("K26").Value = "12/25-01/25"

# Row 30
$ws.Range("A30").Value = 3.36
$ws.Range("B30").Value = 98.3
$ws.Range("C30").Value = -316.55
$ws.Range("D30").Value = 3797.16
$ws.Range("E30").Value = $null
$ws.Range("G30").Value = 1.01
$ws.Range("H30").Value = 13.98
$ws.Range("I30").Value = 18.87
$ws.Range("J30").Value = 3.01
$ws.Range("K30").Value = "12/25-01/25"

# Row 31
$ws.Range("A31").Value = 3.6
$ws.Range("B31").Value = 99.2
$ws.Range("C31").Value = -73.28
$ws.Range("D31").Value = 10852.41
$ws.Range("E31").Value = $null
$ws.Range("G31").Value = 1.03
$ws.Range("H31").Value = 14.03
$ws.Range("I31").Value = 18.73
$ws.Range("J31").Value = 2.94
$ws.Range("K31").Value = "12/25-01/25"

# Row 32
$ws.Range("B32").Value = 95.8
$ws.Range("C32").Value = 378.98
$ws.Range("D32").Value = 2877.38
$ws.Range("E32").Value = $null
$ws.Range("I32").Value = 18.43
$ws.Range("J32").Value = 2.91
$ws.Range("K32").Value = "12/25-01/25"

# Row 33
$ws.Range("B33").Value = 99.5
$ws.Range("C33").Value = -48.83
$ws.Range("D33").Value = 1248.4
$ws.Range("E33").Value = $null
$ws.Range("J33").Value = 2.88
$ws.Range("K33").Value = "12/25-01/25"
$ws.Range("L33").Value = "אלטשולר שחם השתלמות עוקב מדד S&P 500"
$ws.Range("M33").Value = 14862

# Row 34
$ws.Range("A34").Value = 2.28
$ws.Range("B34").Value = 97.8
$ws.Range("C34").Value = -352.39
$ws.Range("D34").Value = 2513.18
$ws.Range("E34").Value = $null
$ws.Range("G34").Value = 0.91
$ws.Range("H34").Value = 12.83
$ws.Range("I34").Value = 18.11
$ws.Range("J34").Value = 2.84
$ws.Range("K34").Value = "12/25-01/25"
$ws.Range("L34").Value = "מור השתלמות -עוקב מדד S&P 500"
$ws.Range("M34").Value = 9451

# Row 35
$ws.Range("C35").Value = 893.87
$ws.Range("D35").Value = 3591.08
$ws.Range("E35").Value = $null
$ws.Range("J35").Value = 2.73
$ws.Range("K35").Value = "12/25-01/25"
$ws.Range("L35").Value = "מגדל השתלמות עוקב מדד S&P 500"
$ws.Range("M35").Value = 14668

# Row 36
$ws.Range("B36").Value = 98.6
$ws.Range("C36").Value = 1155.23
$ws.Range("D36").Value = 7454
$ws.Range("E36").Value = $null
$ws.Range("I36").Value = 18.8
$ws.Range("J36").Value = 2.64
$ws.Range("K36").Value = "12/25-01/25"
$ws.Range("L36").Value = "הראל השתלמות - עוקב מדד s&p 500"
$ws.Range("M36").Value = 13502

# Row 37
$ws.Range("B37").Value = 99.5
$ws.Range("C37").Value = -94.01
$ws.Range("D37").Value = 1735.07
$ws.Range("E37").Value = $null
$ws.Range("I37").Value = 18.13
$ws.Range("J37").Value = 2.55
$ws.Range("K37").Value = "12/25-01/25"

# Row 38
$ws.Range("B38").Value = 99.3
$ws.Range("C38").Value = -95.7
$ws.Range("D38").Value = 1072.35
$ws.Range("E38").Value = $null
$ws.Range("J38").Value = 2.54
$ws.Range("K38").Value = "12/25-01/25"

# Row 39
$ws.Range("A39").Value = 3.12
$ws.Range("B39").Value = 98.6
$ws.Range("C39").Value = 444.5
$ws.Range("D39").Value = 4631.11
$ws.Range("E39").Value = $null
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 13.63
$ws.Range("I39").Value = 18.3
$ws.Range("J39").Value = 2.51
$ws.Range("K39").Value = "12/25-01/25"
$ws.Range("L39").Value = "כלל השתלמות עוקב  מדד s&p 500"
$ws.Range("M39").Value = 13342

# Row 40
$ws.Range("C40").Value = "1,891.8"
$ws.Range("D40").Value = 39772.14
$ws.Range("E40").Value = $null
$ws.Range("H40").Value = "13.98%"
$ws.Range("I40").Value = "18.67%"
$ws.Range("J40").Value = "2.80%"
$ws.Range("K40").Value = "12/25-01/25"
